$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.350.51'
$ws.Range("E2").Value = '  -3.17%  '

$ws.Range("D3").Value = '2.468.17'
$ws.Range("E3").Value = '  -2.18%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.05'
$ws.Range("E5").Value = '  +1.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.65'
$ws.Range("E6").Value = '  -6.77%  '

$ws.Range("E7").Value = '  -3.02%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("E9").Value = '  -4.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.63'
$ws.Range("E10").Value = '  -6.44%  '

$ws.Range("E11").Value = '  -2.69%  '

$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.04'
$ws.Range("E13").Value = '  -3.62%  '

$ws.Range("D14").Value = '2.851.14'
$ws.Range("E14").Value = '  -2.12%  '

$ws.Range("D15").Value = '2.465.35'
$ws.Range("E15").Value = '  -2.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.69'
$ws.Range("E16").Value = '  -6.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.787'
$ws.Range("E17").Value = '  -2.58%  '

$ws.Range("D18").Value = '41.319.52'
$ws.Range("E18").Value = '  -3.16%  '

$ws.Range("E19").Value = '  -5.98%  '

$ws.Range("D20").Value = '0.0₃0922'
$ws.Range("E20").Value = '  -2.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.54'
$ws.Range("E21").Value = '  -5.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.12'
$ws.Range("E22").Value = '  -1.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.84'
$ws.Range("E23").Value = '  -2.41%  '

$ws.Range("E24").Value = '  -3.23%  '

$ws.Range("E25").Value = '  -5.12%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.47'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.24'
$ws.Range("E28").Value = '  -3.87%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.71'
$ws.Range("E29").Value = '  -4.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.17'
$ws.Range("E30").Value = '  -7.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.54'
$ws.Range("E31").Value = '  -2.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.60'
$ws.Range("E32").Value = '  -3.29%  '

$ws.Range("E33").Value = '  -6.38%  '

$ws.Range("E34").Value = '  -0.97%  '

$ws.Range("E35").Value = '  -3.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.03'
$ws.Range("E36").Value = '  -5.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.90'
$ws.Range("E37").Value = '  -6.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '16.95'
$ws.Range("E38").Value = '  -6.59%  '

$ws.Range("E39").Value = '  -6.25%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.114'
$ws.Range("E40").Value = '  -3.69%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.31'
$ws.Range("E41").Value = '  +2.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.34'
$ws.Range("E42").Value = '  -2.76%  '

$ws.Range("E43").Value = '  +0.24%  '

$ws.Range("D44").Value = '1.992.52'
$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("E45").Value = '  -4.40%  '

$ws.Range("E46").Value = '  -5.99%  '

$ws.Range("E47").Value = '  -1.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '70.12'
$ws.Range("E48").Value = '  -2.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '76.09'
$ws.Range("E49").Value = '  -5.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '97.11'
$ws.Range("E50").Value = '  -4.00%  '

$ws.Range("E51").Value = '  -5.85%  '
